$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text (names/URLs) - safe to assign directly.
# Column D holds numeric-looking values that must stay Text (to preserve
# formatting like trailing zeros and double-dot "thousands" notation), so we
# force the cell to Text format before assigning.
# Column E values (e.g. "  +0.73%  ") already stay Text due to the padding.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.213.78'
$ws.Range("E2").Value = '  +0.73%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.602.47'
$ws.Range("E3").Value = '  +0.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.88'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3758'
$ws.Range("E7").Value = '  -0.50%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.84'
$ws.Range("E8").Value = '  +3.72%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3593'
$ws.Range("E9").Value = '  -1.23%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.259'
$ws.Range("E10").Value = '  +0.79%  '

# Row 11
$ws.Range("E11").Value = '  +0.01%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08125'
$ws.Range("E12").Value = '  -0.16%  '

# Row 13
$ws.Range("E13").Value = '  +1.88%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.577'
$ws.Range("E14").Value = '  +0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.332'
$ws.Range("E15").Value = '  -0.42%  '

# Row 16
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.605.84'
$ws.Range("E17").Value = '  +0.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.92'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06918'
$ws.Range("E19").Value = '  +1.05%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.11'
$ws.Range("E20").Value = '  -0.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.519'
$ws.Range("E21").Value = '  +0.27%  '

# Row 22
$ws.Range("E22").Value = '  +0.23%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.86'
$ws.Range("E23").Value = '  -1.21%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.218.94'
$ws.Range("E24").Value = '  +0.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("E25").Value = '  +2.08%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.041'
$ws.Range("E26").Value = '  +9.69%  '

# Row 27
$ws.Range("E27").Value = '  +0.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.51'
$ws.Range("E28").Value = '  +0.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.261'
$ws.Range("E29").Value = '  +0.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.80'
$ws.Range("E30").Value = '  -0.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.407'
$ws.Range("E31").Value = '  +2.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.709'
$ws.Range("E32").Value = '  -0.86%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.779.95'

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9493'
$ws.Range("E34").Value = '  -0.80%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02761'
$ws.Range("E35").Value = '  +2.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07401'
$ws.Range("E36").Value = '  -1.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.25'
$ws.Range("E37").Value = '  +1.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2508'
$ws.Range("E38").Value = '  -0.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.083'
$ws.Range("E39").Value = '  -1.69%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08742'
$ws.Range("E40").Value = '  -0.83%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.399'
$ws.Range("E41").Value = '  +3.19%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7080'
$ws.Range("E42").Value = '  +0.79%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.39'
$ws.Range("E43").Value = '  +0.97%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.84'
$ws.Range("E44").Value = '  +4.63%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6505'
$ws.Range("E45").Value = '  -0.97%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.320'
$ws.Range("E46").Value = '  +2.24%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.008'
$ws.Range("E48").Value = '  +0.29%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.84'
$ws.Range("E49").Value = '  +1.60%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07950'
$ws.Range("E50").Value = '  +0.29%  '

# Row 51
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.192'
$ws.Range("E51").Value = '  -2.36%  '
